$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 607.22
$ws.Range("C3").Value = 635.0700000000001
$ws.Range("C4").Value = 590.08
$ws.Range("C5").Value = 617.96
$ws.Range("C6").Value = 617.96
